$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from "EM-8" to "EF-8,EM-8"
$ws.Range("B9").Value = "EF-8,EM-8"
$ws.Range("C9").Value = "EF-8,EM-8"

# Remove the "Requisitos" section (rows 22-24): "Requisitos:" label and its two entries
$ws.Range("A22:C24").EntireRow.Delete()
